# Apply the commit: "overwrite figures with past ones with right size"
# 1) Shorten two label strings in column A
# 2) Update the B2:B12 numeric figures to their corrected (past) precision values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update shortened labels (column A) ---
$ws.Range("A11").Value = "National tax on millionaires"
$ws.Range("A12").Value = "Global tax on millionaires"

# --- Update numeric figures in column B (rows 2-12) ---
$ws.Range("B2").Value = 0.0460734910149551
$ws.Range("B3").Value = 0.363475325351156
$ws.Range("B4").Value = 0.172116107147519
$ws.Range("B5").Value = -0.12441159299716
$ws.Range("B6").Value = 0.0781277396872659
$ws.Range("B7").Value = 0.168477401346791
$ws.Range("B8").Value = 0.190580455975378
$ws.Range("B9").Value = 0.179302456455172
$ws.Range("B10").Value = -0.118306561742397
$ws.Range("B11").Value = 0.63299052696081
$ws.Range("B12").Value = 0.471390713228356
